$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, pushing existing rows 107-167 down to 108-168
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with the new data record
$ws.Cells.Item(107, 1).Value = 4
$ws.Cells.Item(107, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(107, 3).Value = "Los Lagos"
$ws.Cells.Item(107, 4).Value = 44510
$ws.Cells.Item(107, 5).Value = 10
$ws.Cells.Item(107, 6).Value = "Fruta"
$ws.Cells.Item(107, 7).Value = 100104
$ws.Cells.Item(107, 8).Value = "Frutos de pepita"
$ws.Cells.Item(107, 9).Value = 100104005
$ws.Cells.Item(107, 10).Value = "Pera"
$ws.Cells.Item(107, 11).Value = "Packham's Triumph"
$ws.Cells.Item(107, 12).Value = "Primera"
$ws.Cells.Item(107, 13).Value = 120
$ws.Cells.Item(107, 14).Value = 15000
$ws.Cells.Item(107, 15).Value = 16000
$ws.Cells.Item(107, 16).Value = 15500
$ws.Cells.Item(107, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(107, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(107, 19).Value = 1033
$ws.Cells.Item(107, 20).Value = 15
